$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the authoritative diff of the workbook's XML.
# Values that look like plain numbers (e.g. "211.07") are prefixed with a
# leading apostrophe so Excel stores them as literal text (matching the
# original inline-string cells) instead of silently converting them to
# numeric values.
$updates = @(
    @{ Cell = "D2"; Value = "26.632.22" }
    @{ Cell = "E2"; Value = "  -0.02%  " }
    @{ Cell = "D3"; Value = "1.596.59" }
    @{ Cell = "E3"; Value = "  +0.11%  " }
    @{ Cell = "E4"; Value = "  +0.18%  " }
    @{ Cell = "D5"; Value = "'211.07" }
    @{ Cell = "E5"; Value = "  -0.47%  " }
    @{ Cell = "E6"; Value = "  +0.49%  " }
    @{ Cell = "E7"; Value = "  +0.12%  " }
    @{ Cell = "E8"; Value = "  -0.03%  " }
    @{ Cell = "E9"; Value = "  -0.23%  " }
    @{ Cell = "D10"; Value = "'19.45" }
    @{ Cell = "E10"; Value = "  -1.07%  " }
    @{ Cell = "D11"; Value = "'0.0837" }
    @{ Cell = "E11"; Value = "  +0.12%  " }
    @{ Cell = "D12"; Value = "1.820.73" }
    @{ Cell = "E12"; Value = "  +0.10%  " }
    @{ Cell = "D13"; Value = "1.569.32" }
    @{ Cell = "E13"; Value = "  -1.63%  " }
    @{ Cell = "E14"; Value = "  -0.06%  " }
    @{ Cell = "E15"; Value = "  -0.36%  " }
    @{ Cell = "D16"; Value = "'64.93" }
    @{ Cell = "E16"; Value = "  -0.32%  " }
    @{ Cell = "D17"; Value = "26.617.82" }
    @{ Cell = "E17"; Value = "  -0.02%  " }
    @{ Cell = "E18"; Value = "  +0.65%  " }
    @{ Cell = "E19"; Value = "  +0.24%  " }
    @{ Cell = "D20"; Value = "'208.39" }
    @{ Cell = "E20"; Value = "  -0.68%  " }
    @{ Cell = "D21"; Value = "'7.05" }
    @{ Cell = "E21"; Value = "  +5.22%  " }
    @{ Cell = "D22"; Value = "'4.26" }
    @{ Cell = "E22"; Value = "  +0.32%  " }
    @{ Cell = "D23"; Value = "'2.30" }
    @{ Cell = "E23"; Value = "  -1.16%  " }
    @{ Cell = "D24"; Value = "'8.90" }
    @{ Cell = "E24"; Value = "  +0.16%  " }
    @{ Cell = "D25"; Value = "'145.27" }
    @{ Cell = "E26"; Value = "  +0.10%  " }
    @{ Cell = "D27"; Value = "'7.15" }
    @{ Cell = "E27"; Value = "  -0.11%  " }
    @{ Cell = "E28"; Value = "  +0.19%  " }
    @{ Cell = "E30"; Value = "  +0.31%  " }
    @{ Cell = "E31"; Value = "  +0.16%  " }
    @{ Cell = "E32"; Value = "  -0.43%  " }
    @{ Cell = "E33"; Value = "  +0.76%  " }
    @{ Cell = "D34"; Value = "1.276.34" }
    @{ Cell = "E34"; Value = "  -1.46%  " }
    @{ Cell = "D35"; Value = "'0.617" }
    @{ Cell = "E35"; Value = "  -7.50%  " }
    @{ Cell = "E37"; Value = "  -0.31%  " }
    @{ Cell = "E38"; Value = "  -0.66%  " }
    @{ Cell = "E39"; Value = "  +0.69%  " }
    @{ Cell = "E40"; Value = "  +2.24%  " }
    @{ Cell = "E41"; Value = "  +16.66%  " }
    @{ Cell = "E42"; Value = "  +0.52%  " }
    @{ Cell = "D43"; Value = "'0.785" }
    @{ Cell = "E43"; Value = "  -0.93%  " }
    @{ Cell = "D44"; Value = "'64.14" }
    @{ Cell = "E44"; Value = "  +1.00%  " }
    @{ Cell = "D45"; Value = "1.732.97" }
    @{ Cell = "E45"; Value = "  +0.08%  " }
    @{ Cell = "D46"; Value = "'90.15" }
    @{ Cell = "E46"; Value = "  +0.81%  " }
    @{ Cell = "E47"; Value = "  -0.22%  " }
    @{ Cell = "E48"; Value = "  +3.56%  " }
    @{ Cell = "E49"; Value = "  +1.08%  " }
    @{ Cell = "B50"; Value = "EnergySwap" }
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" }
    @{ Cell = "D50"; Value = "'7.50" }
    @{ Cell = "E50"; Value = "  +0.09%  " }
    @{ Cell = "B51"; Value = "USDD" }
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd" }
    @{ Cell = "D51"; Value = "'1.00" }
    @{ Cell = "E51"; Value = "  +0.21%  " }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
